# ISRaD Milton_1998.xlsx -- "Updated soil type data"
#
# 1. Add a new "pro_usda_soil_order" column to the "profile" sheet
#    (inserted immediately before the existing "pro_soil_taxon" column)
#    and record the USDA soil order ("Spodosols") for the Pickering
#    Nuclear Generating Station profile (row 5).
# 2. Wrap the long citation text in the "metadata" sheet's M4 cell and
#    size the row so the wrapped text is fully visible.

$wb = $excel.ActiveWorkbook

# --- profile sheet -------------------------------------------------------
$wsProfile = $wb.Worksheets.Item("profile")

# Insert a new column before the existing "pro_soil_taxon" column (N)
$wsProfile.Columns("N").Insert()

$wsProfile.Range("N1").Value = "pro_usda_soil_order"
$wsProfile.Range("N5").Value = "Spodosols"

$wsProfile.Select() | Out-Null
$wsProfile.Range("G6").Select() | Out-Null

# --- metadata sheet -------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item("metadata")

$wsMetadata.Range("M4").WrapText = $true
$wsMetadata.Rows.Item(4).RowHeight = 331.2

$wsMetadata.Select() | Out-Null
$wsMetadata.Range("A4").Select() | Out-Null
